# board.pptx edit script
#
# Reproduces the semantically-meaningful, COM-reachable parts of the
# target change:
#   - the empty "Title 25" title placeholder is removed from the slide
#   - the three label rectangles ("Board" / "squares: Square[100]=null" /
#     the empty caption box underneath) are moved as a group to their new
#     spot on the slide and renamed to match the target names
#   - the deck's "first slide number" is reset to the (implicit) default
#
# NOTE: several parts of the upstream diff (removal of the notes master /
# notes slide / now-unused theme2.xml, the en-US -> en-SG locale refresh
# that PowerPoint stamps across every slide layout on resave, new
# p14:creationId GUIDs, etc.) are artifacts of PowerPoint re-serialising
# parts of the package that are not reachable through the documented
# PowerPoint Object Model, so they are intentionally not reproduced here.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Remove the empty title placeholder shape ("Title 25").
$title = $s.Shapes.Item("Title 25")
$title.Delete()

# Move + rename the three remaining rectangles as a group (uniform
# translation of dx=2739433 EMU, dy=-2279630 EMU), matching the target
# offsets exactly (values below are the closest points-precision inputs
# that round-trip to the exact target EMU coordinates).
$board = $s.Shapes.Item("Rectangle 12")
$board.Name = "Rectangle 108"
$board.Left = 251.70339212677166
$board.Top = 201.96094518188977

$squares = $s.Shapes.Item("Rectangle 29")
$squares.Name = "Rectangle 109"
$squares.Left = 251.70339212677166
$squares.Top = 226.8796921393701

$caption = $s.Shapes.Item("Rectangle 4")
$caption.Name = "Rectangle 110"
$caption.Left = 251.70339212677166
$caption.Top = 255.96094518188977

# Reset the deck's first slide number back to its (implicit) default.
$p.PageSetup.FirstSlideNumber = 1
